$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.084.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.795.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "360.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.65%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.233.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.796.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.950"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.998.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0987"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("E30").Value = "  +5.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "51.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.08%  "

$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0851"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "

$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.73%  "

$ws.Range("E40").Value = "  -2.79%  "

$ws.Range("E41").Value = "  +1.92%  "

$ws.Range("E42").Value = "  -1.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.37%  "

$ws.Range("E44").Value = "  -2.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.085.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("E47").Value = "  -2.38%  "

$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.937"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("E51").Value = "  +0.62%  "
